$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.469.34"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +12.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.830.77"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +8.18%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.26"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.546"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.16%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.64"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.97%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.76"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.16%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.284"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.23%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0679"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +8.07%  "

$ws.Range("E12").Value = "  +3.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.088.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.826.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.649"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.90%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "34.458.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +12.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.34"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.37"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +8.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.23"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.89%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "260.85"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0753"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.48%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.40"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.22%  "

$ws.Range("E25").Value = "  +0.59%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.84"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.85"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.55%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.21"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.44%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.117"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.53%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.85"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +10.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0518"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.31%  "

$ws.Range("E33").Value = "  +7.11%  "

$ws.Range("E34").Value = "  +8.44%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.586.63"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.85"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +5.98%  "

$ws.Range("E37").Value = "  +3.54%  "

$ws.Range("E38").Value = "  +5.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.635"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.58"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.07%  "

$ws.Range("E41").Value = "  +5.34%  "

$ws.Range("E42").Value = "  +1.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.921"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +7.48%  "

$ws.Range("E44").Value = "  +6.35%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0521"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.71%  "

$ws.Range("E46").Value = "  +4.27%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.979.28"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +8.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.77"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.97%  "

$ws.Range("E49").Value = "  +1.46%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0124"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +7.14%  "
